# Adds the 6 new LOGT time-log entries (rows 6-11) to the active worksheet,
# matching the target edit described in the diff for
# tspi/ciclo-1/forma-logt1-20106065.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Cells.Item(6, 1).Value = 41899
$ws.Cells.Item(6, 2).Value = 0.340277777777778
$ws.Cells.Item(6, 3).Value = 0.350694444444444
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Formula = "=((HOUR(C6)-HOUR(B6))*60)+(MINUTE(C6)-MINUTE(B6))-D6"
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 8).Value = "Cree la plantilla para las agendas de las reuniones con el cliente."

# Row 7
$ws.Cells.Item(7, 1).Value = 41899
$ws.Cells.Item(7, 2).Value = 0.361111111111111
$ws.Cells.Item(7, 3).Value = 0.395833333333333
$ws.Cells.Item(7, 4).Value = 10
$ws.Cells.Item(7, 5).Formula = "=((HOUR(C7)-HOUR(B7))*60)+(MINUTE(C7)-MINUTE(B7))-D7"
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(7, 8).Value = "Cree la plantilla para las minutas de las reuniones con el cliente."

# Row 8
$ws.Cells.Item(8, 1).Value = 41899
$ws.Cells.Item(8, 2).Value = 0.583333333333333
$ws.Cells.Item(8, 3).Value = 0.600694444444444
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Formula = "=((HOUR(C8)-HOUR(B8))*60)+(MINUTE(C8)-MINUTE(B8))-D8"
$ws.Cells.Item(8, 6).Value = 3
$ws.Cells.Item(8, 8).Value = "Cree la agenda para la reunión #1 con el cliente."

# Row 9
$ws.Cells.Item(9, 1).Value = 41905
$ws.Cells.Item(9, 2).Value = 0.645833333333333
$ws.Cells.Item(9, 3).Value = 0.670138888888889
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Formula = "=((HOUR(C9)-HOUR(B9))*60)+(MINUTE(C9)-MINUTE(B9))-D9"
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(9, 8).Value = "Participé en la reunión #1 con el cliente, y cree la minuta de la reunión #1."

# Row 10
$ws.Cells.Item(10, 1).Value = 41910
$ws.Cells.Item(10, 2).Value = 0.583333333333333
$ws.Cells.Item(10, 3).Value = 0.666666666666667
$ws.Cells.Item(10, 4).Value = 25
$ws.Cells.Item(10, 5).Formula = "=((HOUR(C10)-HOUR(B10))*60)+(MINUTE(C10)-MINUTE(B10))-D10"
$ws.Cells.Item(10, 6).Value = 5
$ws.Cells.Item(10, 8).Value = "Cree el esquema del documento de requerimientos."

# Row 11
$ws.Cells.Item(11, 1).Value = 41910
$ws.Cells.Item(11, 2).Value = 0.666666666666667
$ws.Cells.Item(11, 3).Value = 0.75
$ws.Cells.Item(11, 4).Value = 30
$ws.Cells.Item(11, 5).Formula = "=((HOUR(C11)-HOUR(B11))*60)+(MINUTE(C11)-MINUTE(B11))-D11"
$ws.Cells.Item(11, 6).Value = 6
$ws.Cells.Item(11, 8).Value = "Participé en el analisis de los requerimientos obtenidos en la reunión #1. Se documento una parte del diagrama de casos de uso, y falto documentar los escenarios."

# Row heights to match the wrapped comment text in column H
$ws.Rows.Item(1).RowHeight = 13.75
$ws.Rows.Item(2).RowHeight = 13.75
$ws.Rows.Item(6).RowHeight = 26.65
$ws.Rows.Item(7).RowHeight = 26.65
$ws.Rows.Item(8).RowHeight = 14.15
$ws.Rows.Item(9).RowHeight = 26.65
$ws.Rows.Item(10).RowHeight = 14.15
$ws.Rows.Item(11).RowHeight = 39.15

# Update the selection to reflect where the cursor ended up after data entry
$ws.Range("E12").Select()
